# Laboratorio 7 - Entrega Final
# Insert a new empty paragraph (styled like "Prrafodelista") right after the
# paragraph that asks about the relationship between the number of elements
# in the tree and its height, and right before the existing empty paragraph
# that follows it.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Qué relación encuentra entre el número de elementos en el árbol y la altura del árbol?*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $r.Collapse(0)                 # wdCollapseEnd
    $newRange = $r.InsertParagraphAfter()

    # Re-fetch the newly created paragraph via its range so we can set
    # paragraph-level formatting (style, spacing, alignment) and the
    # run-level formatting used for the character properties stored in
    # the empty paragraph mark.
    $newPara = $newRange
    $newPara.set_Style("Prrafodelista")
    $newPara.ParagraphFormat.SpaceAfter = 0
    $newPara.ParagraphFormat.Alignment = 3   # wdAlignParagraphJustify

    $newPara.Font.Name = "Dax-Regular"
    $newPara.LanguageID = 12298               # es-CO
}
